# Applies the crypto price/volume refresh described in the commit diff.
# Plain numeric-looking Price values are prefixed with a leading apostrophe
# so Excel keeps storing them as text (matching the original inlineStr cells)
# instead of silently re-typing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $text) {
    $ws.Range($address).Value = $text
}

# Row 2
Set-TextValue 'D2' '69.812.52'

# Row 3
Set-TextValue 'D3' '3.496.54'
Set-TextValue 'E3' '  -1.89%  '

# Row 4
Set-TextValue 'E4' '  -0.12%  '

# Row 5
Set-TextValue 'D5' '''602.29'
Set-TextValue 'E5' '  -1.21%  '

# Row 6
Set-TextValue 'D6' '''196.46'
Set-TextValue 'E6' '  +5.59%  '

# Row 7
Set-TextValue 'E7' '  +0.90%  '

# Row 8
Set-TextValue 'E8' '  -0.09%  '

# Row 9
Set-TextValue 'D9' '''0.209'
Set-TextValue 'E9' '  -2.74%  '

# Row 10
Set-TextValue 'D10' '''0.652'
Set-TextValue 'E10' '  +1.28%  '

# Row 11
Set-TextValue 'D11' '''54.05'
Set-TextValue 'E11' '  +0.44%  '

# Row 12
Set-TextValue 'D12' '''0.0000301'
Set-TextValue 'E12' '  -2.59%  '

# Row 13
Set-TextValue 'D13' '''9.56'
Set-TextValue 'E13' '  +0.82%  '

# Row 14
Set-TextValue 'D14' '4.055.66'
Set-TextValue 'E14' '  -1.80%  '

# Row 15
Set-TextValue 'D15' '''601.25'
Set-TextValue 'E15' '  +3.67%  '

# Row 16
Set-TextValue 'D16' '69.940.33'
Set-TextValue 'E16' '  -0.46%  '

# Row 17
Set-TextValue 'D17' '''19.08'
Set-TextValue 'E17' '  +0.57%  '

# Row 18
Set-TextValue 'D18' '''12.59'
Set-TextValue 'E18' '  -0.75%  '

# Row 19
Set-TextValue 'D19' '3.505.43'
Set-TextValue 'E19' '  -1.54%  '

# Row 20
Set-TextValue 'E20' '  +0.23%  '

# Row 21
Set-TextValue 'D21' '''0.990'
Set-TextValue 'E21' '  -0.41%  '

# Row 22
Set-TextValue 'D22' '''18.24'
Set-TextValue 'E22' '  +5.11%  '

# Row 23
Set-TextValue 'D23' '''104.76'
Set-TextValue 'E23' '  +10.40%  '

# Row 24
Set-TextValue 'B24' 'Toncoin'
Set-TextValue 'C24' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D24' '''5.03'
Set-TextValue 'E24' '  +3.49%  '

# Row 25
Set-TextValue 'B25' 'PancakeSwap'
Set-TextValue 'C25' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D25' '''4.58'
Set-TextValue 'E25' '  -2.87%  '

# Row 26
Set-TextValue 'D26' '''3.10'
Set-TextValue 'E26' '  +5.25%  '

# Row 27
Set-TextValue 'D27' '''10.96'
Set-TextValue 'E27' '  -0.10%  '

# Row 28
Set-TextValue 'E28' '  +3.79%  '

# Row 29
Set-TextValue 'D29' '''33.58'
Set-TextValue 'E29' '  +4.01%  '

# Row 30
Set-TextValue 'D30' '''4.46'
Set-TextValue 'E30' '  +21.25%  '

# Row 31
Set-TextValue 'D31' '''7.20'
Set-TextValue 'E31' '  +2.44%  '

# Row 32
Set-TextValue 'D32' '''12.70'
Set-TextValue 'E32' '  +3.69%  '

# Row 33
Set-TextValue 'E33' '  +0.83%  '

# Row 34
Set-TextValue 'D34' '''63.63'
Set-TextValue 'E34' '  +0.17%  '

# Row 35
Set-TextValue 'D35' '3.734.62'
Set-TextValue 'E35' '  +5.62%  '

# Row 36
Set-TextValue 'B36' 'Dai'
Set-TextValue 'C36' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D36' '''0.999'
Set-TextValue 'E36' '  -0.25%  '

# Row 37
Set-TextValue 'B37' 'PEPE'
Set-TextValue 'C37' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D37' '0.0₃0802'
Set-TextValue 'E37' '  +2.33%  '

# Row 38
Set-TextValue 'D38' '''509.86'
Set-TextValue 'E38' '  -4.14%  '

# Row 39
Set-TextValue 'E39' '  -7.76%  '

# Row 40
Set-TextValue 'D40' '''0.390'
Set-TextValue 'E40' '  -3.27%  '

# Row 41
Set-TextValue 'D41' '''36.59'
Set-TextValue 'E41' '  -1.43%  '

# Row 42
Set-TextValue 'D42' '''3.51'
Set-TextValue 'E42' '  -0.38%  '

# Row 43
Set-TextValue 'D43' '''0.136'
Set-TextValue 'E43' '  +0.50%  '

# Row 44
Set-TextValue 'D44' '''0.0456'
Set-TextValue 'E44' '  -0.82%  '

# Row 45
Set-TextValue 'B45' 'Stellar'
Set-TextValue 'C45' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D45' '''0.140'
Set-TextValue 'E45' '  -0.54%  '

# Row 46
Set-TextValue 'B46' 'ApeXProtocol'
Set-TextValue 'C46' 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 'D46' '''3.31'
Set-TextValue 'E46' '  -4.63%  '

# Row 47
Set-TextValue 'B47' 'ThetaToken'
Set-TextValue 'C47' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D47' '''2.81'
Set-TextValue 'E47' '  -3.72%  '

# Row 48
Set-TextValue 'E48' '  +0.32%  '

# Row 49
Set-TextValue 'D49' '''8.73'

# Row 50
Set-TextValue 'D50' '''131.93'
Set-TextValue 'E50' '  -3.29%  '

# Row 51
Set-TextValue 'D51' '''0.000240'
Set-TextValue 'E51' '  -2.14%  '

Write-Output "Updated $([int]104) cells across $((($ws.UsedRange.Rows.Count)))-row sheet"
